# Natmi following Dr Hou advice
# Update the Col1a1-Itga11 LR-pairs sheet: recompute existing rows 2-4 with the
# revised statistics, and add new rows 5-7 to cover the full combination of
# sending/target clusters (ECs/FAPs/sCs) x (FAPs/sCs) target values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col1a1"
$ws.Range("C2").Value = "Itga11"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6.072131
$ws.Range("H2").Value = 18.216393
$ws.Range("I2").Value = 0.003943999267036455
$ws.Range("J2").Value = 0.003943999267036454
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 24.759128
$ws.Range("N2").Value = 74.277384
$ws.Range("O2").Value = 0.9895671066967037
$ws.Range("P2").Value = 0.9895671066967037
$ws.Range("Q2").Value = 150.340668661768
$ws.Range("R2").Value = 1353.066017955912
$ws.Range("S2").Value = 0.003902851943495184
$ws.Range("T2").Value = 0.003902851943495184

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col1a1"
$ws.Range("C3").Value = "Itga11"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.072131
$ws.Range("H3").Value = 18.216393
$ws.Range("I3").Value = 0.003943999267036455
$ws.Range("J3").Value = 0.003943999267036454
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2610326666666667
$ws.Range("N3").Value = 0.7830980000000001
$ws.Range("O3").Value = 0.0104328933032964
$ws.Range("P3").Value = 0.0104328933032964
$ws.Range("Q3").Value = 1.585024547279333
$ws.Range("R3").Value = 14.265220925514
$ws.Range("S3").Value = [double]"4.114732354127055e-05"
$ws.Range("T3").Value = [double]"4.114732354127054e-05"

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Col1a1"
$ws.Range("C4").Value = "Itga11"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1480.851806666667
$ws.Range("H4").Value = 4442.55542
$ws.Range("I4").Value = 0.9618498744646554
$ws.Range("J4").Value = 0.9618498744646552
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 24.759128
$ws.Range("N4").Value = 74.277384
$ws.Range("O4").Value = 0.9895671066967037
$ws.Range("P4").Value = 0.9895671066967037
$ws.Range("Q4").Value = 36664.59943029125
$ws.Range("R4").Value = 329981.3948726212
$ws.Range("S4").Value = 0.9518149973505767
$ws.Range("T4").Value = 0.9518149973505765

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col1a1"
$ws.Range("C5").Value = "Itga11"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1480.851806666667
$ws.Range("H5").Value = 4442.55542
$ws.Range("I5").Value = 0.9618498744646554
$ws.Range("J5").Value = 0.9618498744646552
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2610326666666667
$ws.Range("N5").Value = 0.7830980000000001
$ws.Range("O5").Value = 0.0104328933032964
$ws.Range("P5").Value = 0.0104328933032964
$ws.Range("Q5").Value = 386.5506960323511
$ws.Range("R5").Value = 3478.95626429116
$ws.Range("S5").Value = 0.01003487711407879
$ws.Range("T5").Value = 0.01003487711407879

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Col1a1"
$ws.Range("C6").Value = "Itga11"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 52.663316
$ws.Range("H6").Value = 157.989948
$ws.Range("I6").Value = 0.03420612626830831
$ws.Range("J6").Value = 0.0342061262683083
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 24.759128
$ws.Range("N6").Value = 74.277384
$ws.Range("O6").Value = 0.9895671066967037
$ws.Range("P6").Value = 0.9895671066967037
$ws.Range("Q6").Value = 1303.897781748448
$ws.Range("R6").Value = 11735.08003573603
$ws.Range("S6").Value = 0.03384925740263197
$ws.Range("T6").Value = 0.03384925740263196

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Col1a1"
$ws.Range("C7").Value = "Itga11"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 52.663316
$ws.Range("H7").Value = 157.989948
$ws.Range("I7").Value = 0.03420612626830831
$ws.Range("J7").Value = 0.0342061262683083
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.2610326666666667
$ws.Range("N7").Value = 0.7830980000000001
$ws.Range("O7").Value = 0.0104328933032964
$ws.Range("P7").Value = 0.0104328933032964
$ws.Range("Q7").Value = 13.74684581098933
$ws.Range("R7").Value = 123.721612298904
$ws.Range("S7").Value = 0.0003568688656763449
$ws.Range("T7").Value = 0.0003568688656763448
